$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Five new "Form responses 1" survey rows (56-60) appended at the bottom of the sheet.
# Each entry lists the column letter -> cell value for that row, in sheet column order.
$newRows = [ordered]@{
    56 = [ordered]@{
        A = 45059.87175856481
        B = 'Raghavi '
        C = 'raghavi_c@yahoo.com'
        D = 'Yes'
        E = 'Yes'
        F = 'Yes'
        H = 'Maybe'
        I = 'Coupons to obtain other external products/services'
        J = 'Every time you purchase some product'
        K = 1.0
        L = 'No'
        M = 'Very Important'
        N = 3.0
        O = 'Yes'
        P = 'Very important'
        Q = 'Yes'
        R = 'Expected to invest very high to avail loyalty program '
        S = 'Moderate Understanding'
        T = 4.0
        U = 'The rewarding scheme should be smooth irrespective of background technology used'
        V = 'No'
        W = 'Yes'
        X = 'Maybe'
        Y = 'Maybe'
    }
    57 = [ordered]@{
        A = 45061.50101116898
        B = 'Sreedeva Krupananda B Reddy'
        C = '20bcs128@iiitdwd.ac.in'
        D = 'Yes'
        E = 'Yes'
        F = 'Yes'
        H = 'Yes'
        I = 'Cash Incentive, Discount on the next purchases/service'
        J = 'Every time you purchase some product'
        K = 3.0
        L = 'No'
        M = 'Somewhat important'
        N = 3.0
        O = 'No'
        P = 'Very important'
        Q = 'No'
        S = 'Moderate Understanding'
        T = 3.0
        U = 'Transparency, Security and increased Flexibility'
        V = 'No'
        W = 'Yes'
        X = 'Yes'
        Y = 'No'
    }
    58 = [ordered]@{
        A = 45064.33640958334
        B = 'Aalekh Prasad'
        C = '20bec001@iiitdwd.ac.in'
        D = 'Yes'
        E = 'No'
        H = 'No'
        I = 'Discount on the next purchases/service'
        J = 'Monthly'
        K = 2.0
        L = 'No'
        M = 'Very Important'
        N = 3.0
        O = 'Maybe'
        P = 'Very important'
        Q = 'No'
        S = 'Moderate Understanding'
        T = 4.0
        U = 'discount on next transaction'
        V = 'No'
        W = 'Yes'
        X = 'No'
        Y = 'Yes'
    }
    59 = [ordered]@{
        A = 45064.99178840278
        B = 'JAGADISH'
        C = '20CS032@iiitdwd.ac.in'
        D = 'Yes'
        E = 'Yes'
        F = 'Yes'
        H = 'Yes'
        I = 'Cash Incentive'
        J = 'Every time you purchase some product'
        K = 4.0
        L = 'Yes'
        M = 'Very Important'
        N = 4.0
        O = 'Yes'
        P = 'Very important'
        Q = 'Yes'
        R = 'meeting the limitations '
        S = 'Moderate Understanding'
        T = 4.0
        U = 'Faster Redemption and Lower Fees'
        V = 'Yes'
        W = 'Yes'
        X = 'Yes'
        Y = 'Yes'
    }
    60 = [ordered]@{
        A = 45092.84549443287
        B = 'Anuj'
        C = '19bcs012@iiitdwd.ac.in'
        D = 'Yes'
        E = 'Yes'
        F = 'Yes'
        H = 'Maybe'
        I = 'Cash Incentive'
        J = 'Every time you purchase some product'
        K = 3.0
        L = 'No'
        M = 'Not important'
        N = 1.0
        O = 'Maybe'
        P = 'Very important'
        Q = 'Yes'
        R = 'wallet issues as in the reward can be used for particular wallet or type of payment method.'
        S = 'Moderate Understanding'
        T = 4.0
        U = 'product authenticity.'
        V = 'Yes'
        W = 'Maybe'
        X = 'No'
        Y = 'Maybe'
    }
}

# Reference cells whose formatting (number format / font) the new rows should copy:
#  - column A uses the Timestamp date/time style
#  - every other populated column uses the shared general text style
$dateFormatSource = $ws.Range("A2")
$textFormatSource = $ws.Range("B2")

foreach ($rowNum in $newRows.Keys) {
    $rowCells = $newRows[$rowNum]
    foreach ($col in $rowCells.Keys) {
        $ws.Range("$col$rowNum").Value = $rowCells[$col]
    }
}

foreach ($rowNum in $newRows.Keys) {
    $rowCells = $newRows[$rowNum]
    foreach ($col in $rowCells.Keys) {
        if ($col -eq "A") {
            $dateFormatSource.Copy()
        } else {
            $textFormatSource.Copy()
        }
        $ws.Range("$col$rowNum").PasteSpecial(-4122)
    }
}

$excel.CutCopyMode = $false